# Add 13 new GPU rows (RX 7801 XT .. RX 7813 XT) to the "SPEC" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SPEC")

for ($i = 0; $i -lt 13; $i++) {
    $row = 4 + $i
    $ws.Cells.Item($row, 1).Value = "RX " + (7801 + $i) + " XT"
    $ws.Cells.Item($row, 2).Value = 111 + $i
    $ws.Cells.Item($row, 3).Value = 91 + $i
    $ws.Cells.Item($row, 4).Value = 56 + $i
}

# Match the author's final selection state (single cell G10 selected).
$ws.Range("G10").Select()
